$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the BOM line item for H12/H13 "WireHole_PTH" (row 12) - no longer
# need to buy this part, so the whole row is deleted and subsequent rows
# shift up.
$ws.Rows.Item(12).Select() | Out-Null
$ws.Rows.Item(12).Delete()

# Row deletion does not automatically shrink the conditional formatting
# range that covered the deleted row, so fix it up explicitly (E2:E40 -> E2:E39).
$fc = $ws.Range("E2:E40").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("E2:E39"))
